$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32 (hunk 1)
$ws.Range("H32").Value = 1396
$ws.Range("I32").Value = 800
$ws.Range("J32").Value = 1992
$ws.Range("K32").Value = 800
$ws.Range("L32").Value = 1992
$ws.Range("M32").Value = -474
$ws.Range("N32").Value = -2644
# Row 40 (hunk 2)
$ws.Range("H40").Value = 3600
$ws.Range("I40").Value = 3600
$ws.Range("J40").Value = 3600
$ws.Range("K40").Value = 3600
$ws.Range("L40").Value = 3600
$ws.Range("M40").Value = -3425
$ws.Range("N40").Value = -3950
# Row 80 (hunk 3)
$ws.Range("H80").Value = 397.3889
$ws.Range("I80").Value = 305.75
$ws.Range("J80").Value = 580.6667
$ws.Range("K80").Value = 917.25
$ws.Range("L80").Value = 1742.0001
$ws.Range("M80").Value = 80.75
$ws.Range("N80").Value = -3738.0001
# Row 83 (hunk 4)
$ws.Range("H83").Value = 397.3889
$ws.Range("I83").Value = 305.75
$ws.Range("J83").Value = 580.6667
$ws.Range("K83").Value = 2751.75
$ws.Range("L83").Value = 5226.0003
$ws.Range("M83").Value = 2240.25
$ws.Range("N83").Value = -15210.0003
# Row 98 (hunk 5)
$ws.Range("H98").Value = 705.65
$ws.Range("I98").Value = 701.94116
$ws.Range("J98").Value = 726.6667
$ws.Range("K98").Value = 701.94116
$ws.Range("L98").Value = 726.6667
$ws.Range("M98").Value = 796.05884
$ws.Range("N98").Value = -3722.6667
# Row 122 (hunk 6)
$ws.Range("H122").Value = 705.65
$ws.Range("I122").Value = 701.94116
$ws.Range("J122").Value = 726.6667
$ws.Range("K122").Value = 2105.82348
$ws.Range("L122").Value = 2180.0001
$ws.Range("M122").Value = 344.17652
$ws.Range("N122").Value = -7080.0001
# Row 127 (hunk 7)
$ws.Range("H127").Value = 1374.7273
$ws.Range("I127").Value = 558
$ws.Range("J127").Value = 2055.3333
$ws.Range("K127").Value = 1674
$ws.Range("L127").Value = 6165.999899999999
$ws.Range("M127").Value = 3286
$ws.Range("N127").Value = -16085.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (hunk 8)
$ws.Range("H61").Value = 111334940
$ws.Range("I61").Value = 166834080
$ws.Range("K61").Value = 166834080
$ws.Range("M61").Value = -166833868
# Row 136 (hunk 9)
$ws.Range("H136").Value = 111334940
$ws.Range("I136").Value = 166834080
$ws.Range("K136").Value = 500502240
$ws.Range("M136").Value = -500499690

$ws = $wb.Worksheets.Item("CRP")
# Row 9 (hunk 10)
$ws.Range("H9").Value = 36548.668
$ws.Range("J9").Value = 36548.668
$ws.Range("L9").Value = 36548.668
$ws.Range("N9").Value = -36884.668
# Row 16 (hunk 11)
$ws.Range("H16").Value = 1314
$ws.Range("I16").Value = 506.5
$ws.Range("K16").Value = 506.5
$ws.Range("M16").Value = -219.5
# Row 22 (hunk 12)
$ws.Range("H22").Value = 90909730
$ws.Range("I22").Value = 200000510
$ws.Range("J22").Value = 733
$ws.Range("K22").Value = 200000510
$ws.Range("L22").Value = 733
$ws.Range("M22").Value = -200000160
$ws.Range("N22").Value = -1433
# Row 105 (hunk 13)
$ws.Range("H105").Value = 741.6842
$ws.Range("I105").Value = 785.4666999999999
$ws.Range("J105").Value = 577.5
$ws.Range("K105").Value = 785.4666999999999
$ws.Range("L105").Value = 577.5
$ws.Range("M105").Value = 961.5333000000001
$ws.Range("N105").Value = -4071.5
# Row 106 (hunk 14)
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").ClearContents()
$ws.Range("N106").Value = 0
# Row 107 (hunk 15)
$ws.Range("H107").Value = 366.66666
$ws.Range("I107").Value = 335.7143
$ws.Range("J107").Value = 438.8889
$ws.Range("K107").Value = 335.7143
$ws.Range("L107").Value = 438.8889
$ws.Range("M107").Value = 1584.2857
$ws.Range("N107").Value = -4278.8889
# Row 113 (hunk 16)
$ws.Range("H113").Value = 1314
$ws.Range("I113").Value = 506.5
$ws.Range("K113").Value = 506.5
$ws.Range("M113").Value = 1663.5
# Row 122 (hunk 17)
$ws.Range("H122").Value = 1890.5834
$ws.Range("I122").Value = 1396.625
$ws.Range("J122").Value = 2878.5
$ws.Range("K122").Value = 4189.875
$ws.Range("L122").Value = 8635.5
$ws.Range("M122").Value = -1739.875
$ws.Range("N122").Value = -13535.5
# Row 132 (hunk 18)
$ws.Range("H132").Value = 37264.586
$ws.Range("I132").Value = 2420.0435
$ws.Range("K132").Value = 7260.130500000001
$ws.Range("M132").Value = -4730.130500000001

$ws = $wb.Worksheets.Item("GSM")
# Row 102 (hunk 19)
$ws.Range("H102").Value = 868.6
$ws.Range("I102").Value = 766.2727
$ws.Range("K102").Value = 766.2727
$ws.Range("M102").Value = 855.7273
# Row 126 (hunk 20)
$ws.Range("H126").Value = 1951.4
$ws.Range("I126").Value = 1785.7142
$ws.Range("J126").Value = 2338
$ws.Range("K126").Value = 5357.142599999999
$ws.Range("L126").Value = 7014
$ws.Range("M126").Value = -2887.142599999999
$ws.Range("N126").Value = -11954

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (hunk 21)
$ws.Range("H46").Value = 1000
$ws.Range("J46").Value = 1000
$ws.Range("L46").Value = 1000
$ws.Range("N46").Value = -1376

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (hunk 22)
$ws.Range("H81").Value = 1816.8276
$ws.Range("I81").Value = 686.55554
$ws.Range("J81").Value = 2325.45
$ws.Range("K81").Value = 1373.11108
$ws.Range("L81").Value = 4650.9
$ws.Range("M81").Value = -312.1110799999999
$ws.Range("N81").Value = -6772.9
# Row 84 (hunk 23)
$ws.Range("H84").Value = 1816.8276
$ws.Range("I84").Value = 686.55554
$ws.Range("J84").Value = 2325.45
$ws.Range("K84").Value = 6865.555399999999
$ws.Range("L84").Value = 23254.5
$ws.Range("M84").Value = -1561.555399999999
$ws.Range("N84").Value = -33862.5
# Row 107 (hunk 24)
$ws.Range("H107").Value = 289.1111
$ws.Range("I107").Value = 251.11111
$ws.Range("J107").Value = 327.1111
$ws.Range("K107").Value = 753.3333299999999
$ws.Range("L107").Value = 981.3333
$ws.Range("M107").Value = 1166.66667
$ws.Range("N107").Value = -4821.3333
# Row 125 (hunk 25)
$ws.Range("H125").Value = 46905
$ws.Range("J125").Value = 46905
$ws.Range("L125").Value = 46905
$ws.Range("N125").Value = -56745
# Row 126 (hunk 26)
$ws.Range("H126").Value = 1042.0952
$ws.Range("I126").Value = 868.38464
$ws.Range("J126").Value = 1324.375
$ws.Range("K126").Value = 2605.15392
$ws.Range("L126").Value = 3973.125
$ws.Range("M126").Value = -135.1539199999997
$ws.Range("N126").Value = -8913.125
# Row 136 (hunk 27)
$ws.Range("H136").Value = 44725.195
$ws.Range("I136").Value = 25322.781
$ws.Range("K136").Value = 75968.34299999999
$ws.Range("M136").Value = -73418.34299999999
# Row 138 (hunk 28)
$ws.Range("H138").Value = 41166.668
$ws.Range("J138").Value = 41166.668
$ws.Range("L138").Value = 41166.668
$ws.Range("N138").Value = -51446.668
# Row 139 (hunk 29)
$ws.Range("H139").Value = 29750
$ws.Range("J139").Value = 29750
$ws.Range("L139").Value = 29750
$ws.Range("N139").Value = -40030
# Row 140 (hunk 30)
$ws.Range("H140").Value = 56530
$ws.Range("J140").Value = 56530
$ws.Range("L140").Value = 56530
$ws.Range("N140").Value = -66890
# Row 141 (hunk 31)
$ws.Range("H141").Value = 50416.668
$ws.Range("J141").Value = 50416.668
$ws.Range("L141").Value = 50416.668
$ws.Range("N141").Value = -60776.668
